$d = $word.ActiveDocument
$styles = $d.Styles

# The document-wide East Asian font used throughout the stylesheet changes
# from "DejaVu Sans" to "Tahoma". Apply it to every paragraph style whose
# rPr currently carries an explicit eastAsia font (Normal, Heading).
$eastAsiaFont = "Tahoma"

$normal = $styles.Item("Normal")
$normal.Font.NameFarEast = $eastAsiaFont

$heading = $styles.Item("Heading")
$heading.Font.NameFarEast = $eastAsiaFont

# List, Caption and Index previously had no explicit rFonts override; they
# now pick up an explicit complex-script (cs) font matching the existing
# "DejaVu Sans" value (inherited/default) so it becomes explicit in rPr.
$csFont = "DejaVu Sans"

$list = $styles.Item("List")
$list.Font.NameBi = $csFont

$caption = $styles.Item("Caption")
$caption.Font.NameBi = $csFont

$index = $styles.Item("Index")
$index.Font.NameBi = $csFont
